# Natmi following Dr Hou advice
# Update the LR-pair edge statistics on the active sheet (Ccl25-Ccr9) to reflect
# the recomputed values (Ligand/Receptor-expressing cell counts changed from 1 to 3,
# with the corresponding average/total expression, specificity and edge-weight figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; E=3; G=5.934604333333334; H=17.803813; I=0.3081877218757661; J=0.3081877218757661; K=3; M=1.381632333333333; N=4.144897; O=0.4379959344090866; P=0.4379959344090866; Q=8.199441232473445; R=73.79497109226101; S=0.1349849692163839; T=0.1349849692163839 }
    @{ Row=3; E=3; G=5.934604333333334; H=17.803813; I=0.3081877218757661; J=0.3081877218757661; K=3; M=1.136233333333333; N=3.4087; O=0.3602011682365698; P=0.3602011682365697; Q=6.743095263677778; R=60.6878573731; S=0.111009577455818; T=0.111009577455818 }
    @{ Row=4; E=3; G=5.934604333333334; H=17.803813; I=0.3081877218757661; J=0.3081877218757661; K=3; M=0.6365753333333333; N=1.909726; O=0.2018028973543437; P=0.2018028973543437; Q=3.777822731693111; R=34.000404585238; S=0.06219317520356426; T=0.06219317520356425 }
    @{ Row=5; E=3; G=7.720664; H=23.161992; I=0.4009389195777736; J=0.4009389195777736; K=3; M=1.381632333333333; N=4.144897; O=0.4379959344090866; P=0.4379959344090866; Q=10.66711901720267; R=96.00407115482402; S=0.1756096167214366; T=0.1756096167214366 }
    @{ Row=6; E=3; G=7.720664; H=23.161992; I=0.4009389195777736; J=0.4009389195777736; K=3; M=1.136233333333333; N=3.4087; O=0.3602011682365698; P=0.3602011682365697; Q=8.772475792266667; R=78.9522821304; S=0.1444186672234221; T=0.1444186672234221 }
    @{ Row=7; E=3; G=7.720664; H=23.161992; I=0.4009389195777736; J=0.4009389195777736; K=3; M=0.6365753333333333; N=1.909726; O=0.2018028973543437; P=0.2018028973543437; Q=4.914784259354667; R=44.233058334192; S=0.08091063563291491; T=0.0809106356329149 }
    @{ Row=8; E=3; G=5.601191; H=16.803573; I=0.2908733585464604; J=0.2908733585464603; K=3; M=1.381632333333333; N=4.144897; O=0.4379959344090866; P=0.4379959344090866; Q=7.738786590775667; R=69.649079316981; S=0.1274013484712662; T=0.1274013484712662 }
    @{ Row=9; E=3; G=5.601191; H=16.803573; I=0.2908733585464604; J=0.2908733585464603; K=3; M=1.136233333333333; N=3.4087; O=0.3602011682365698; P=0.3602011682365697; Q=6.364259920566666; R=57.2783392851; S=0.1047729235573297; T=0.1047729235573296 }
    @{ Row=10; E=3; G=5.601191; H=16.803573; I=0.2908733585464604; J=0.2908733585464603; K=3; M=0.6365753333333333; N=1.909726; O=0.2018028973543437; P=0.2018028973543437; Q=3.565580027888666; R=32.09022025099799; S=0.05869908651786457; T=0.05869908651786455 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("E" + $r).Value = $u.E
    $ws.Range("G" + $r).Value = $u.G
    $ws.Range("H" + $r).Value = $u.H
    $ws.Range("I" + $r).Value = $u.I
    $ws.Range("J" + $r).Value = $u.J
    $ws.Range("K" + $r).Value = $u.K
    $ws.Range("M" + $r).Value = $u.M
    $ws.Range("N" + $r).Value = $u.N
    $ws.Range("O" + $r).Value = $u.O
    $ws.Range("P" + $r).Value = $u.P
    $ws.Range("Q" + $r).Value = $u.Q
    $ws.Range("R" + $r).Value = $u.R
    $ws.Range("S" + $r).Value = $u.S
    $ws.Range("T" + $r).Value = $u.T
}
